# Refresh the crypto price/volume snapshot (GitHub Actions bot).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2, E2
$ws.Range("D2").Value = "30.118.61"
$ws.Range("E2").Value = "  -0.72%  "

# Row 3: D3, E3
$ws.Range("D3").Value = "1.908.70"
$ws.Range("E3").Value = "  -1.42%  "

# Row 4: D4, E4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.00%  "

# Row 5: D5, E5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7414"
$ws.Range("E5").Value = "  -1.40%  "

# Row 6: D6, E6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "244.24"
$ws.Range("E6").Value = "  +0.39%  "

# Row 7: E7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8: D8, E8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3088"
$ws.Range("E8").Value = "  -3.37%  "

# Row 9: D9, E9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.52"
$ws.Range("E9").Value = "  -5.04%  "

# Row 10: D10, E10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06977"
$ws.Range("E10").Value = "  -0.97%  "

# Row 11: D11, E11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08080"
$ws.Range("E11").Value = "  +0.45%  "

# Row 12: D12, E12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7698"
$ws.Range("E12").Value = "  -1.76%  "

# Row 13: D13, E13
$ws.Range("D13").Value = "1.940.31"
$ws.Range("E13").Value = "  +0.32%  "

# Row 14: D14, E14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.322"
$ws.Range("E14").Value = "  -1.72%  "

# Row 15: D15, E15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.33"
$ws.Range("E15").Value = "  -0.98%  "

# Row 16: D16, E16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.27"
$ws.Range("E16").Value = "  -1.58%  "

# Row 17: D17, E17
$ws.Range("D17").Value = "30.126.48"
$ws.Range("E17").Value = "  -0.69%  "

# Row 18: D18, E18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.111"
$ws.Range("E18").Value = "  +0.42%  "

# Row 19: D19, E19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007828"
$ws.Range("E19").Value = "  -2.50%  "

# Row 20: D20, E20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.01"
$ws.Range("E20").Value = "  -5.08%  "

# Row 21: D21, E21
$ws.Range("D21").Value = "2.131.83"
$ws.Range("E21").Value = "  -2.26%  "

# Row 22: E22
$ws.Range("E22").Value = "  +0.09%  "

# Row 23: D23, E23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.04%  "

# Row 24: D24, E24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.133"
$ws.Range("E24").Value = "  +6.24%  "

# Row 25: D25, E25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.386"
$ws.Range("E25").Value = "  -1.86%  "

# Row 26: D26, E26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.99"
$ws.Range("E26").Value = "  +1.25%  "

# Row 27: D27, E27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.97"
$ws.Range("E27").Value = "  -0.80%  "

# Row 28: D28, E28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1274"
$ws.Range("E28").Value = "  -2.44%  "

# Row 29: D29, E29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.045"
$ws.Range("E29").Value = "  -7.82%  "

# Row 30: D30, E30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.548"
$ws.Range("E30").Value = "  +0.87%  "

# Row 31: D31, E31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.352"
$ws.Range("E31").Value = "  -1.73%  "

# Row 32: D32, E32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.337"
$ws.Range("E32").Value = "  -2.96%  "

# Row 33: D33, E33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.082"
$ws.Range("E33").Value = "  -1.76%  "

# Row 34: D34, E34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.308"
$ws.Range("E34").Value = "  -2.43%  "

# Row 35: D35, E35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05158"
$ws.Range("E35").Value = "  -2.13%  "

# Row 36: D36, E36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7482"
$ws.Range("E36").Value = "  -1.47%  "

# Row 37: D37, E37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.724"
$ws.Range("E37").Value = "  -2.21%  "

# Row 38: D38, E38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01958"
$ws.Range("E38").Value = "  -0.29%  "

# Row 39: D39, E39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.798"
$ws.Range("E39").Value = "  -0.40%  "

# Row 40: D40, E40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.346"
$ws.Range("E40").Value = "  -3.56%  "

# Row 41: D41, E41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4503"
$ws.Range("E41").Value = "  -0.53%  "

# Row 42: D42, E42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.51"
$ws.Range("E42").Value = "  -5.82%  "

# Row 43: E43
$ws.Range("E43").Value = "  -0.39%  "

# Row 44: D44, E44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.08%  "

# Row 45: D45, E45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8395"
$ws.Range("E45").Value = "  +0.25%  "

# Row 46: D46, E46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.741"
$ws.Range("E46").Value = "  +0.06%  "

# Row 47: B47, C47, D47, E47
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.88"
$ws.Range("E47").Value = "  +0.13%  "

# Row 48: B48, C48, D48, E48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.926"
$ws.Range("E48").Value = "  -1.00%  "

# Row 49: D49, E49
$ws.Range("D49").Value = "2.074.44"
$ws.Range("E49").Value = "  -0.36%  "

# Row 50: E50
$ws.Range("E50").Value = "  -2.68%  "

# Row 51: D51, E51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1187"
$ws.Range("E51").Value = "  -3.93%  "

